$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.623.30'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').Value = '2.360.65'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '522.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.539'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  +5.24%  '
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.36'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').Value = '2.783.34'
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('D15').Value = '57.595.85'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').Value = '2.366.96'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '331.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.08%  '
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('E24').Value = '  +4.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.991'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.91%  '
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('E30').Value = '  +0.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  +3.29%  '
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +4.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '38.70'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '151.56'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.384'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '284.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.30%  '
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0221'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.50%  '
$ws.Range('E51').Value = '  +0.00%  '
